$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (row 2) ---
$ws.Range('I2').Value = 'is even?'
$ws.Range('J2').Value = 'repeat'
$ws.Range('K2').Value = 'log_cube'

# --- New formula columns (rows 3-12) ---
for ($r = 3; $r -le 12; $r++) {
    $ws.Range("I$r").Formula = '=IF(MOD(B' + $r + ',2)=0,"yes","no")'
    $ws.Range("J$r").Formula = '=REPT("|",G' + $r + ')'
    $ws.Range("K$r").Formula = '=LOG10(E' + $r + ')'
}

# --- Column J width (auto-fit like the original author did, ~26.57 chars wide) ---
$null = $ws.Columns.Item(10).AutoFit()
$ws.Columns.Item(10).ColumnWidth = 25.72

# --- New defined names ---
$wb.Names.Add('is_even', '=Sheet1!$I$3')
$wb.Names.Add('repeat_line', '=Sheet1!$J$3')
$wb.Names.Add('log_cube', '=Sheet1!$K$3')

# --- Selection moves to J3, matching the saved view state ---
$null = $ws.Range('J3').Select()
